# Update cryptocurrency price/volume data in the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where both Price (column D) and Volume(1h) (column E) change
$priceAndVolumeUpdates = @(
    @{Row=2; D="29.898.57"; E="  -0.15%  "},
    @{Row=3; D="1.875.52"; E="  -0.69%  "},
    @{Row=4; D="1.000"; E="  -0.13%  "},
    @{Row=5; D="0.7391"; E="  -4.32%  "},
    @{Row=6; D="242.30"; E="  -0.21%  "},
    @{Row=7; D="0.9997"; E="  -0.18%  "},
    @{Row=8; D="0.3160"; E="  +2.16%  "},
    @{Row=9; D="0.07184"; E="  +0.90%  "},
    @{Row=10; D="24.72"; E="  -3.31%  "},
    @{Row=11; D="0.08373"; E="  -2.20%  "},
    @{Row=12; D="0.7511"; E="  -1.52%  "},
    @{Row=13; D="1.935.07"; E="  -3.81%  "},
    @{Row=14; D="5.424"; E="  +1.87%  "},
    @{Row=15; D="92.66"; E="  -1.05%  "},
    @{Row=16; D="29.903.89"; E="  -0.62%  "},
    @{Row=17; D="6.078"; E="  -1.22%  "},
    @{Row=18; D="247.16"; E="  +1.39%  "},
    @{Row=19; D="13.57"; E="  -1.16%  "},
    @{Row=20; D="0.000007829"; E="  +0.86%  "},
    @{Row=21; D="0.9990"; E="  +0.04%  "},
    @{Row=22; D="2.126.69"; E="  -6.18%  "},
    @{Row=23; D="8.010"; E="  +0.65%  "},
    @{Row=24; D="1.000"; E="  -0.16%  "},
    @{Row=25; D="0.1549"; E="  -5.78%  "},
    @{Row=26; D="9.270"; E="  -0.31%  "},
    @{Row=27; D="165.10"; E="  +1.67%  "},
    @{Row=28; D="18.66"; E="  -0.38%  "},
    @{Row=29; D="2.037"; E="  +0.43%  "},
    @{Row=30; D="1.502"; E="  +4.88%  "},
    @{Row=31; D="4.594"; E="  +2.16%  "},
    @{Row=32; D="1.534"; E="  -0.35%  "},
    @{Row=33; D="4.285"; E="  +4.33%  "},
    @{Row=34; D="0.05316"; E="  -2.12%  "},
    @{Row=35; D="1.238"; E="  -0.16%  "},
    @{Row=36; D="0.7551"; E="  +1.20%  "},
    @{Row=37; D="1.001"; E="  -0.50%  "},
    @{Row=39; D="0.01962"; E="  +0.32%  "},
    @{Row=40; D="2.757"; E="  -1.05%  "},
    @{Row=41; D="0.4503"; E="  +1.00%  "},
    @{Row=42; D="1.111.30"; E="  +0.44%  "},
    @{Row=43; D="6.062"; E="  -0.26%  "},
    @{Row=44; D="72.43"; E="  -1.22%  "},
    @{Row=45; D="0.8532"; E="  +0.68%  "},
    @{Row=47; D="103.09"; E="  -0.32%  "},
    @{Row=48; D="7.622"; E="  +0.13%  "},
    @{Row=49; D="1.846"; E="  -1.20%  "},
    @{Row=50; D="2.978"; E="  -0.76%  "}
)

foreach ($item in $priceAndVolumeUpdates) {
    $r = $item.Row
    $dCell = $ws.Cells.Item($r, 4)
    # Force column D to be treated as literal text so values like
    # "1.000" or "242.30" are not silently converted to numbers
    # and lose their trailing zeros / formatting.
    $dCell.NumberFormat = "@"
    $dCell.Value = $item.D
    $dCell.ClearFormats()
    $ws.Cells.Item($r, 5).Value = $item.E
}

# Rows where only Volume(1h) (column E) changes
$volumeOnlyUpdates = @(
    @{Row=38; E="  -0.08%  "},
    @{Row=46; E="  -0.06%  "}
)

foreach ($item in $volumeOnlyUpdates) {
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}

# Row 51: coin replaced (RocketPoolETH -> EnergySwap) with new link, price and volume
$ws.Cells.Item(51, 2).Value = "EnergySwap"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$d51 = $ws.Cells.Item(51, 4)
$d51.NumberFormat = "@"
$d51.Value = "9.488"
$d51.ClearFormats()
$ws.Cells.Item(51, 5).Value = "  -2.60%  "

